$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "64.672.70"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.761.62"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'577.09"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'160.36"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "3.251.82"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'27.36"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "64.264.68"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "2.768.83"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'12.18"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "'4.86"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "'358.66"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -7.19%  "
$ws.Range("D24").Value = "'65.17"
$ws.Range("D25").Value = "'0.172"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'8.63"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "0.0₃0929"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "'7.39"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.38"
$ws.Range("E30").Value = "  +9.15%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.97"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'167.68"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "'5.02"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").Value = "'20.23"
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'1.86"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "'352.77"
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").Value = "'6.44"
$ws.Range("E40").Value = "  +4.37%  "
$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").Value = "'39.16"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'22.67"
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").Value = "'21.64"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "'0.0594"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "'136.57"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'0.632"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "2.145.23"
$ws.Range("E51").Value = "  +0.65%  "
